# Renamed Thermdat to Nasa. Changed excel import function. Bug fixes.
# Concretely: simplify the long fully-qualified Python class names used for
# the thermodynamic model in column F down to short friendly names, and
# update the description of the "thermo_model" column (F2) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-5 used the fully qualified IdealGasThermo class name -> "IdealGas"
$ws.Range("F3").Value = "IdealGas"
$ws.Range("F4").Value = "IdealGas"
$ws.Range("F5").Value = "IdealGas"

# Rows 6-11 used the fully qualified HarmonicThermo class name -> "Harmonic"
$ws.Range("F6").Value = "Harmonic"
$ws.Range("F7").Value = "Harmonic"
$ws.Range("F8").Value = "Harmonic"
$ws.Range("F9").Value = "Harmonic"
$ws.Range("F10").Value = "Harmonic"
$ws.Range("F11").Value = "Harmonic"

# Update the header description in F2 (row 2 holds human readable
# descriptions of each column).
$ws.Range("F2").Value = "Type of thermodynamic model. Supported options include IdealGas and Harmonic"

# Reflect the final selection/view state left by the editor: cell F2 selected
# and the view scrolled back to show column A (no special top-left cell).
$null = $ws.Range("F2").Select()
